$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column M: URL built from DOI in column A
$ws.Range("M2").Formula = '=(CONCATENATE("https://www.jacc.org/doi/",A2))'
$ws.Range("M3:M11").Formula = '=(CONCATENATE("https://www.jacc.org/doi/",A3))'

# Column M width (54.83203125 in raw OOXML units ~= 54 "characters" of ColumnWidth)
$ws.Columns.Item(13).ColumnWidth = 54

# Update selection / view state
$ws.Range("M2:M11").Select()
$excel.ActiveWindow.ScrollColumn = 3
